# Update the "想去人数" (number of people interested) counts (column F)
# on the "展览" sheet and the "全部类型" sheet to reflect newly generated
# output data.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 115
$ws1.Range("F4").Value = 1587
$ws1.Range("F5").Value = 286
$ws1.Range("F6").Value = 79
$ws1.Range("F7").Value = 1903
$ws1.Range("F8").Value = 10265
$ws1.Range("F10").Value = 140
$ws1.Range("F12").Value = 196
$ws1.Range("F14").Value = 7127
$ws1.Range("F16").Value = 671
$ws1.Range("F17").Value = 74
$ws1.Range("F19").Value = 244

# --- Sheet "全部类型" (All Types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 115
$ws4.Range("F4").Value = 1587
$ws4.Range("F5").Value = 286
$ws4.Range("F7").Value = 79
$ws4.Range("F8").Value = 1903
$ws4.Range("F11").Value = 10265
$ws4.Range("F13").Value = 140
$ws4.Range("F15").Value = 196
$ws4.Range("F17").Value = 7127
$ws4.Range("F19").Value = 671
$ws4.Range("F20").Value = 74
$ws4.Range("F22").Value = 244
